$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("K2").Value = 0.00282676936157604
$ws.Range("L2").Value = 0.1254160836934189
$ws.Range("M2").Value = 0.03195169266094546
$ws.Range("N2").Value = 0.1167113203949155
$ws.Range("O2").Value = 0.01289363475846455
$ws.Range("P2").Value = 0.0000911771095807754
$ws.Range("Q2").Value = 0.007980949400260331
$ws.Range("R2").Value = 0.03883417507905353
$ws.Range("S2").Value = 0.06256198248808779
$ws.Range("T2").Value = 0.05510141632740324
$ws.Range("U2").Value = 0.03261024947901932
$ws.Range("V2").Value = 0.07040611201028658
$ws.Range("W2").Value = 0.003861385362598371
$ws.Range("X2").Value = 0.08298421775308297
$ws.Range("Y2").Value = 0.001381502715995762
$ws.Range("Z2").Value = 0.08867687461820604
$ws.Range("AA2").Value = 0.03031628336290294
$ws.Range("AB2").Value = 0.003722121473973689
$ws.Range("AC2").Value = 0.000163118590079386
$ws.Range("AD2").Value = 0.0002928347196976026
$ws.Range("AE2").Value = 0.02259854206519518
$ws.Range("AF2").Value = 0.001305594047853933
$ws.Range("AG2").Value = 0.1055053691779202
$ws.Range("AH2").Value = 0.02316014190598597
$ws.Range("AI2").Value = 0.01594285617016413
$ws.Range("AJ2").Value = 0.008211137021179494
$ws.Range("AK2").Value = 0.01421853412833693
$ws.Range("AL2").Value = 0.00141338090808862
$ws.Range("AM2").Value = 0.02964132368621568
$ws.Range("AN2").Value = 0.00325589770968777
$ws.Range("AO2").Value = 0.005939199696859528
$ws.Range("AP2").Value = 0.0000241221229638384
$ws.Range("AQ2").Value = 0
$ws.Range("K4").Value = 0.003336327932820915
$ws.Range("L4").Value = 0.09953649852449062
$ws.Range("M4").Value = 0.0160213917014956
$ws.Range("N4").Value = 0.1287604787181179
$ws.Range("O4").Value = 0.01870592104309638
$ws.Range("P4").Value = 0.004866687224992905
$ws.Range("Q4").Value = 0.009146257248529206
$ws.Range("R4").Value = 0.01985390025843351
$ws.Range("S4").Value = 0.01018752612964069
$ws.Range("T4").Value = 0.08036073726445821
$ws.Range("U4").Value = 0.006363775302593837
$ws.Range("V4").Value = 0.1051139343853866
$ws.Range("W4").Value = 0.01373603853851814
$ws.Range("X4").Value = 0.01829748249514992
$ws.Range("Y4").Value = 0.0228756567878444
$ws.Range("Z4").Value = 0.04978618167382518
$ws.Range("AA4").Value = 0.07991965684944877
$ws.Range("AB4").Value = 0.0001552516099089872
$ws.Range("AC4").Value = 0.005965587206817133
$ws.Range("AD4").Value = 0.01213651591511327
$ws.Range("AE4").Value = 0.003781075549226944
$ws.Range("AF4").Value = 0.002663739532638173
$ws.Range("AG4").Value = 0.09862425677090191
$ws.Range("AH4").Value = 0.07241950410239698
$ws.Range("AI4").Value = 0.007874900805778233
$ws.Range("AJ4").Value = 0.0210388512807249
$ws.Range("AK4").Value = 0.01535293899276911
$ws.Range("AL4").Value = 0.0007831222056222252
$ws.Range("AM4").Value = 0.05094993159730611
$ws.Range("AN4").Value = 0.01221487004281544
$ws.Range("AO4").Value = 0.00916713053495828
$ws.Range("AP4").Value = 0.000003871774179492559
$ws.Range("AQ4").Value = 0
$ws.Range("K6").Value = 0.00149228814728463
$ws.Range("L6").Value = 0.09351253485406169
$ws.Range("M6").Value = 0.0318994432421267
$ws.Range("N6").Value = 0.1234346366041957
$ws.Range("O6").Value = 0.004848202601762009
$ws.Range("P6").Value = 0.002728318117251663
$ws.Range("Q6").Value = 0.006129443178039842
$ws.Range("R6").Value = 0.01390860024107324
$ws.Range("S6").Value = 0.008790014954384744
$ws.Range("T6").Value = 0.06641389846789152
$ws.Range("U6").Value = 0.008395601828900628
$ws.Range("V6").Value = 0.1229579800054655
$ws.Range("W6").Value = 0.01214236946923035
$ws.Range("X6").Value = 0.05047306027414464
$ws.Range("Y6").Value = 0.02676149311398067
$ws.Range("Z6").Value = 0.06301718668125389
$ws.Range("AA6").Value = 0.06518387408592695
$ws.Range("AB6").Value = 0.00355525584640397
$ws.Range("AC6").Value = 0.00462639847744024
$ws.Range("AD6").Value = 0.01154725075118409
$ws.Range("AE6").Value = 0.007293222914675762
$ws.Range("AF6").Value = 0.002846101420602571
$ws.Range("AG6").Value = 0.0972121266509208
$ws.Range("AH6").Value = 0.06804575741788714
$ws.Range("AI6").Value = 0.009960585388581288
$ws.Range("AJ6").Value = 0.0184779229405993
$ws.Range("AK6").Value = 0.01351005450429638
$ws.Range("AL6").Value = 0.00003691130707428152
$ws.Range("AM6").Value = 0.04609837428074648
$ws.Range("AN6").Value = 0.006782081658755164
$ws.Range("AO6").Value = 0.007872433721655884
$ws.Range("AP6").Value = 0.00004657685220227744
$ws.Range("AQ6").Value = 0

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("K2").Value = 0.00282676936157604
$ws.Range("L2").Value = 0.1282428530549949
$ws.Range("M2").Value = 0.1601945457159404
$ws.Range("N2").Value = 0.2769058661108559
$ws.Range("O2").Value = 0.2897995008693204
$ws.Range("P2").Value = 0.2898906779789012
$ws.Range("Q2").Value = 0.2978716273791616
$ws.Range("R2").Value = 0.3367058024582151
$ws.Range("S2").Value = 0.3992677849463029
$ws.Range("T2").Value = 0.4543692012737061
$ws.Range("U2").Value = 0.4869794507527254
$ws.Range("V2").Value = 0.557385562763012
$ws.Range("W2").Value = 0.5612469481256104
$ws.Range("X2").Value = 0.6442311658786933
$ws.Range("Y2").Value = 0.6456126685946891
$ws.Range("Z2").Value = 0.7342895432128952
$ws.Range("AA2").Value = 0.7646058265757981
$ws.Range("AB2").Value = 0.7683279480497718
$ws.Range("AC2").Value = 0.7684910666398511
$ws.Range("AD2").Value = 0.7687839013595488
$ws.Range("AE2").Value = 0.791382443424744
$ws.Range("AF2").Value = 0.7926880374725979
$ws.Range("AG2").Value = 0.8981934066505181
$ws.Range("AH2").Value = 0.921353548556504
$ws.Range("AI2").Value = 0.9372964047266681
$ws.Range("AJ2").Value = 0.9455075417478476
$ws.Range("AK2").Value = 0.9597260758761845
$ws.Range("AL2").Value = 0.9611394567842731
$ws.Range("AM2").Value = 0.9907807804704888
$ws.Range("AN2").Value = 0.9940366781801766
$ws.Range("AO2").Value = 0.9999758778770361
$ws.Range("AP2").Value = 1
$ws.Range("K4").Value = 0.003336327932820915
$ws.Range("L4").Value = 0.1028728264573115
$ws.Range("M4").Value = 0.1188942181588071
$ws.Range("N4").Value = 0.247654696876925
$ws.Range("O4").Value = 0.2663606179200214
$ws.Range("P4").Value = 0.2712273051450143
$ws.Range("Q4").Value = 0.2803735623935435
$ws.Range("R4").Value = 0.300227462651977
$ws.Range("S4").Value = 0.3104149887816177
$ws.Range("T4").Value = 0.3907757260460759
$ws.Range("U4").Value = 0.3971395013486698
$ws.Range("V4").Value = 0.5022534357340565
$ws.Range("W4").Value = 0.5159894742725746
$ws.Range("X4").Value = 0.5342869567677245
$ws.Range("Y4").Value = 0.557162613555569
$ws.Range("Z4").Value = 0.6069487952293942
$ws.Range("AA4").Value = 0.686868452078843
$ws.Range("AB4").Value = 0.687023703688752
$ws.Range("AC4").Value = 0.6929892908955692
$ws.Range("AD4").Value = 0.7051258068106824
$ws.Range("AE4").Value = 0.7089068823599093
$ws.Range("AF4").Value = 0.7115706218925475
$ws.Range("AG4").Value = 0.8101948786634494
$ws.Range("AH4").Value = 0.8826143827658464
$ws.Range("AI4").Value = 0.8904892835716246
$ws.Range("AJ4").Value = 0.9115281348523495
$ws.Range("AK4").Value = 0.9268810738451186
$ws.Range("AL4").Value = 0.9276641960507409
$ws.Range("AM4").Value = 0.978614127648047
$ws.Range("AN4").Value = 0.9908289976908624
$ws.Range("AO4").Value = 0.9999961282258207
$ws.Range("AP4").Value = 1
$ws.Range("K6").Value = 0.00149228814728463
$ws.Range("L6").Value = 0.09500482300134633
$ws.Range("M6").Value = 0.126904266243473
$ws.Range("N6").Value = 0.2503389028476687
$ws.Range("O6").Value = 0.2551871054494307
$ws.Range("P6").Value = 0.2579154235666823
$ws.Range("Q6").Value = 0.2640448667447222
$ws.Range("R6").Value = 0.2779534669857954
$ws.Range("S6").Value = 0.2867434819401802
$ws.Range("T6").Value = 0.3531573804080717
$ws.Range("U6").Value = 0.3615529822369724
$ws.Range("V6").Value = 0.4845109622424378
$ws.Range("W6").Value = 0.4966533317116682
$ws.Range("X6").Value = 0.5471263919858128
$ws.Range("Y6").Value = 0.5738878850997935
$ws.Range("Z6").Value = 0.6369050717810474
$ws.Range("AA6").Value = 0.7020889458669743
$ws.Range("AB6").Value = 0.7056442017133783
$ws.Range("AC6").Value = 0.7102706001908186
$ws.Range("AD6").Value = 0.7218178509420027
$ws.Range("AE6").Value = 0.7291110738566784
$ws.Range("AF6").Value = 0.731957175277281
$ws.Range("AG6").Value = 0.8291693019282018
$ws.Range("AH6").Value = 0.8972150593460889
$ws.Range("AI6").Value = 0.9071756447346702
$ws.Range("AJ6").Value = 0.9256535676752694
$ws.Range("AK6").Value = 0.9391636221795658
$ws.Range("AL6").Value = 0.9392005334866401
$ws.Range("AM6").Value = 0.9852989077673866
$ws.Range("AN6").Value = 0.9920809894261418
$ws.Range("AO6").Value = 0.9999534231477977
$ws.Range("AP6").Value = 1

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.557385562763012
$ws.Range("D4").Value = 21
$ws.Range("F4").Value = 0.5022534357340565
$ws.Range("G4").Value = 12
$ws.Range("F6").Value = 0.5471263919858128

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.7342895432128952
$ws.Range("D4").Value = 29
$ws.Range("F4").Value = 0.7051258068106824
$ws.Range("G4").Value = 20
$ws.Range("D6").Value = 26
$ws.Range("F6").Value = 0.7020889458669743
$ws.Range("G6").Value = 17

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.8981934066505181
$ws.Range("D4").Value = 32
$ws.Range("F4").Value = 0.8101948786634494
$ws.Range("G4").Value = 23
$ws.Range("F6").Value = 0.8291693019282018

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.921353548556504
$ws.Range("D4").Value = 35
$ws.Range("F4").Value = 0.9115281348523495
$ws.Range("G4").Value = 26
$ws.Range("D6").Value = 34
$ws.Range("F6").Value = 0.9071756447346702
$ws.Range("G6").Value = 25
